$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder the four algorithm rows (2-5): the row that used to hold "Custom
# Algorithm" (old row 4) moves to row 2, "Genetic Algorithm" (old row 5)
# moves to row 3, "Bee Colony Optimization" (old row 2) moves to row 4 and
# "Ant Colony Optimization" (old row 3) moves to row 5 - now carrying real
# measured values instead of the old 1/1/1 placeholders. Each row's
# formatting (the bold "Time" column highlight) travels with its data.

# Row 2 -> Custom Algorithm (formatting swaps with old row 4's)
$ws.Range("A2").Value = "Custom Algorithm"
$ws.Range("B2").Value = 293.81508892713703
$ws.Range("C2").Value = 25.449044320137599
$ws.Range("D2").Value = 738.99285714285702
$ws.Range("G2").Value = "Custom Algorithm"
$ws.Range("H2").Formula = "=B2*D2"
$ws.Range("B2").Font.Bold = $true
$ws.Range("C2").Font.Bold = $true
$ws.Range("D2").Font.Bold = $false
$ws.Range("H2").Font.Bold = $true

# Row 3 -> Genetic Algorithm (formatting unchanged)
$ws.Range("A3").Value = "Genetic Algorithm"
$ws.Range("B3").Value = 1206.7301450085999
$ws.Range("C3").Value = 61.813138841078597
$ws.Range("D3").Value = 777.28214285714296
$ws.Range("G3").Value = "Genetic Algorithm"
$ws.Range("H3").Formula = "=B3*D3"

# Row 4 -> Bee Colony Optimization (formatting swaps with old row 2's)
$ws.Range("A4").Value = "Bee Colony Optimization"
$ws.Range("B4").Value = 5752.8784876649397
$ws.Range("C4").Value = 28.5210469018932
$ws.Range("D4").Value = 726.55357142857099
$ws.Range("G4").Value = "Bee Colony Optimization"
$ws.Range("H4").Formula = "=B4*D4"
$ws.Range("B4").Font.Bold = $false
$ws.Range("C4").Font.Bold = $false
$ws.Range("D4").Font.Bold = $true
$ws.Range("H4").Font.Bold = $false

# Row 5 -> Ant Colony Optimization, now with real computed figures
# (formatting unchanged)
$ws.Range("A5").Value = "Ant Colony Optimization"
$ws.Range("B5").Value = 6540.4494637690204
$ws.Range("C5").Value = 53.399727705482498
$ws.Range("D5").Value = 792.92280314099298
$ws.Range("G5").Value = "Ant Colony Optimization"
$ws.Range("H5").Formula = "=B5*D5"

# Move the active selection, matching the updated sheet view
$ws.Range("G10").Select()
